# 2016 Stronghold Summer Robot Controller Layout - update controls
$wb = $excel.ActiveWorkbook

# --- Sheet "X-Box Controller" ---
$ws1 = $wb.Worksheets.Item("X-Box Controller")

# "Arm Up?" / "Arm Down?" -> "Arm Up" / "Arm Down" (strip the trailing "?")
$ws1.Range("B15").Value = "Arm Up"
$ws1.Range("C15").Value = "Arm Up"
$ws1.Range("B16").Value = "Arm Down"
$ws1.Range("C16").Value = "Arm Down"

# New buttons documented: Ball Intake (row 19) and Fire Ball Shooter (row 20)
$ws1.Range("B19").Value = "Ball Intake"
$ws1.Range("C19").Value = "Ball Intake"
$ws1.Range("B20").Value = "Fire Ball Shooter"
$ws1.Range("C20").Value = "Fire Ball Shooter"

# Cursor / selection on this sheet now rests on C12, and the top-left freeze
# scroll position resets to the sheet origin.
$ws1.Activate() | Out-Null
$ws1.Range("C12").Select() | Out-Null

# --- Sheet "Joysticks" ---
$ws2 = $wb.Worksheets.Item("Joysticks")

# Row 9 ("Arcade" row, button 2) used to say "Arm Down?" in column B only;
# now it is documented per-joystick in columns C and D as "Arm Down".
$ws2.Range("B9").Value = ""
$ws2.Range("C9").ClearFormats()
$ws2.Range("C9").Value = "Arm Down"
$ws2.Range("D9").Value = "Arm Down"

# Row 10 ("Left Joystick" row, button 3) used to say "Arm Up?" in column B
# only; now it is documented per-joystick in columns C and D as "Arm Up".
$ws2.Range("B10").Value = ""
$ws2.Range("C10").Value = "Arm Up"
$ws2.Range("D10").Value = "Arm Up"

# Button 1 (row 8) now documents "Fire Ball Shooter" in columns C and D.
$ws2.Range("C8").Value = "Fire Ball Shooter"
$ws2.Range("D8").Value = "Fire Ball Shooter"

# Button 6 (row 13) now documents "Ball Intake" in columns C and D.
$ws2.Range("C13").Value = "Ball Intake"
$ws2.Range("D13").Value = "Ball Intake"

# The "Joysticks" tab becomes the active/selected tab, with the cursor at D14.
$ws2.Activate() | Out-Null
$ws2.Range("D14").Select() | Out-Null
